# "remove delay when insert/remove safekey"
#
# The audio-clip table on Sheet1 gets two new rows appended for the
# "safekey removed" / "safekey inserted" clips that already have their
# text (column C, rows 31/32) in the sheet but were still missing an
# ID (column A) and an audio file name (column B). This adds those,
# matching the numbering/style pattern used by the rest of the table,
# and refreshes the wrapped-text row heights the same way Excel does
# whenever the table content changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: ID 30 / file "030" (uses the safekey-removed text already in C31) ---
$ws.Range("A30").Copy()
$ws.Range("A31").PasteSpecial(-4122)   # xlPasteFormats - reuse the ID-column style
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "030"

# --- Row 32: ID 31 / file "031" (uses the safekey-inserted text already in C32) ---
$ws.Range("A30").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "031"

$excel.CutCopyMode = $false

# --- Recompute the wrapped-text row heights (Excel re-measures these whenever
#     the sheet is edited/resaved; the values below match the new measurement) ---
$ws.Rows.Item(2).RowHeight = 60
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 45
$ws.Rows.Item(20).RowHeight = 30
$ws.Rows.Item(28).RowHeight = 30
$ws.Rows.Item(29).RowHeight = 30
$ws.Rows.Item(32).RowHeight = 90

# --- Selection moved one column over (D32) after the new columns were filled in ---
$ws.Range("D32").Select()
